# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F3").Value  = 671
$sheet1.Range("F7").Value  = 602
$sheet1.Range("F8").Value  = 51
$sheet1.Range("F9").Value  = 50
$sheet1.Range("F11").Value = 102
$sheet1.Range("F12").Value = 7
$sheet1.Range("F13").Value = 94
$sheet1.Range("F14").Value = 297
$sheet1.Range("F15").Value = 409
$sheet1.Range("F16").Value = 491
$sheet1.Range("F17").Value = 127
$sheet1.Range("F18").Value = 11079
$sheet1.Range("F19").Value = 5316

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F3").Value  = 671
$sheet4.Range("F7").Value  = 602
$sheet4.Range("F8").Value  = 51
$sheet4.Range("F9").Value  = 50
$sheet4.Range("F13").Value = 102
$sheet4.Range("F14").Value = 7
$sheet4.Range("F15").Value = 94
$sheet4.Range("F16").Value = 297
$sheet4.Range("F17").Value = 409
$sheet4.Range("F18").Value = 491
$sheet4.Range("F19").Value = 127
$sheet4.Range("F20").Value = 11079
$sheet4.Range("F22").Value = 5316

$wb.Save()
